$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update D2, E2
$ws.Range('D2').Value = '60.532.63'
$ws.Range('E2').Value = '  -2.21%  '

# Row 3: update D3, E3
$ws.Range('D3').Value = '2.899.34'
$ws.Range('E3').Value = '  -3.27%  '

# Row 4: update E4
$ws.Range('E4').Value = '  +0.11%  '

# Row 5: update D5, E5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '527.46'
$ws.Range('E5').Value = '  -4.15%  '

# Row 6: update D6, E6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.08'
$ws.Range('E6').Value = '  -6.80%  '

# Row 7: update E7
$ws.Range('E7').Value = '  -0.03%  '

# Row 8: update D8, E8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.554'
$ws.Range('E8').Value = '  -2.19%  '

# Row 9: update D9, E9
$ws.Range('D9').Value = '2.906.36'
$ws.Range('E9').Value = '  -3.03%  '

# Row 10: update D10, E10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.108'
$ws.Range('E10').Value = '  -3.68%  '

# Row 11: update D11, E11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.92'
$ws.Range('E11').Value = '  -5.00%  '

# Row 12: update D12, E12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.359'
$ws.Range('E12').Value = '  -1.96%  '

# Row 13: update D13, E13
$ws.Range('D13').Value = '3.405.22'
$ws.Range('E13').Value = '  -3.01%  '

# Row 14: update E14
$ws.Range('E14').Value = '  +1.49%  '

# Row 15: update D15, E15
$ws.Range('D15').Value = '60.547.78'
$ws.Range('E15').Value = '  -2.27%  '

# Row 16: update D16, E16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '22.60'
$ws.Range('E16').Value = '  -4.44%  '

# Row 17: update D17, E17
$ws.Range('D17').Value = '2.894.35'
$ws.Range('E17').Value = '  -3.06%  '

# Row 18: update E18
$ws.Range('E18').Value = '  -4.49%  '

# Row 19: update D19, E19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.95'
$ws.Range('E19').Value = '  -2.96%  '

# Row 20: update D20, E20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.62'
$ws.Range('E20').Value = '  -2.68%  '

# Row 21: update D21, E21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '363.71'
$ws.Range('E21').Value = '  -7.45%  '

# Row 22: update D22, E22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.60'
$ws.Range('E22').Value = '  -0.40%  '

# Row 23: update E23
$ws.Range('E23').Value = '  -0.19%  '

# Row 24: update D24, E24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.44'
$ws.Range('E24').Value = '  -2.37%  '

# Row 25: update D25, E25
$ws.Range('D25').Value = '3.006.70'
$ws.Range('E25').Value = '  -3.67%  '

# Row 26: update D26, E26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.451'
$ws.Range('E26').Value = '  -3.32%  '

# Row 27: update E27
$ws.Range('E27').Value = '  -1.53%  '

# Row 28: update D28, E28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  -0.24%  '

# Row 29: update D29, E29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.86'
$ws.Range('E29').Value = '  -6.94%  '

# Row 30: update D30, E30
$ws.Range('D30').Value = '0.0₃0860'
$ws.Range('E30').Value = '  -9.00%  '

# Row 31: update E31
$ws.Range('E31').Value = '  -0.01%  '

# Row 32: update E32
$ws.Range('E32').Value = '  -2.77%  '

# Row 33: update D33, E33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.52'
$ws.Range('E33').Value = '  -4.53%  '

# Row 34: update D34, E34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '148.27'
$ws.Range('E34').Value = '  -6.86%  '

# Row 35: update D35, E35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.35'
$ws.Range('E35').Value = '  -6.58%  '

# Row 36: update D36, E36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.58'
$ws.Range('E36').Value = '  -7.01%  '

# Row 37: update D37, E37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  -7.25%  '

# Row 38: update E38
$ws.Range('E38').Value = '  -6.23%  '

# Row 39: update D39, E39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '37.88'
$ws.Range('E39').Value = '  +1.50%  '

# Row 40: update D40, E40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.50'
$ws.Range('E40').Value = '  -4.24%  '

# Row 41: update D41, E41
$ws.Range('D41').Value = '2.334.77'
$ws.Range('E41').Value = '  -4.45%  '

# Row 42: update D42, E42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.68'
$ws.Range('E42').Value = '  -5.99%  '

# Row 43: update D43, E43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.644'
$ws.Range('E43').Value = '  -2.43%  '

# Row 44: update D44, E44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '20.73'
$ws.Range('E44').Value = '  -7.52%  '

# Row 45: update D45, E45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0573'
$ws.Range('E45').Value = '  -3.27%  '

# Row 46: update B46, C46, D46, E46
$ws.Range('B46').Value = 'FirstDigitalUSD'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.998'
$ws.Range('E46').Value = '  +0.11%  '

# Row 47: update B47, C47, D47, E47
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.07'
$ws.Range('E47').Value = '  +2.67%  '

# Row 48: update D48, E48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0235'
$ws.Range('E48').Value = '  -4.71%  '

# Row 49: update B49, C49, D49, E49
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0932'
$ws.Range('E49').Value = '  -1.50%  '

# Row 50: update B50, C50, D50, E50
$ws.Range('B50').Value = 'WhiteBITCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '10.33'
$ws.Range('E50').Value = '  -1.48%  '

# Row 51: update D51, E51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '250.24'
$ws.Range('E51').Value = '  -5.67%  '
